$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "date"
$ws.Range("B1").Value = "name"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "x"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "y"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "z"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "a"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "b"

$ws.Range("B6").Select()
